# The diff shows a new weekly price-report row being inserted before the
# existing row 61 (Damasco / Castle Brite / Primera, Región de O'Higgins),
# pushing the former rows 61-75 down to 62-76 unchanged, and extending the
# sheet dimension from A1:T75 to A1:T76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 61, shifting rows 61:75 down to 62:76.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new record's data.
$ws.Cells.Item(61, 1).Value = 4
$ws.Cells.Item(61, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(61, 3).Value = "Los Lagos"
$ws.Cells.Item(61, 4).Value = 45275
$ws.Cells.Item(61, 5).Value = 10
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100103
$ws.Cells.Item(61, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(61, 9).Value = 100103003
$ws.Cells.Item(61, 10).Value = "Damasco"
$ws.Cells.Item(61, 11).Value = "Castle Brite"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 350
$ws.Cells.Item(61, 14).Value = 22000
$ws.Cells.Item(61, 15).Value = 22000
$ws.Cells.Item(61, 16).Value = 22000
$ws.Cells.Item(61, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(61, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(61, 19).Value = 1467
$ws.Cells.Item(61, 20).Value = 15
